$wb = $excel.ActiveWorkbook

# Localize the builtin "Comma" cell style name to "Millares" (Spanish)
try {
    $commaStyle = $wb.Styles.Item("Comma")
    $commaStyle.Name = "Millares"
} catch {
}

# Rename sheet "Funciones 1" -> "Sumas maximos promedios" (2nd sheet, rId2)
$wsFunciones1 = $wb.Worksheets.Item(2)
$wsFunciones1.Name = "Sumas maximos promedios"

# Sheet 1 "Inversión": becomes the active/selected tab, selection moves to K12
$wsInversion = $wb.Worksheets.Item(1)
$wsInversion.Activate()
$wsInversion.Range("K12").Select()

# Sheet 2 "Sumas maximos promedios": selection moves to F32
$wsFunciones1.Range("F32").Select()

# Sheet 3 "Funciones 2": selection moves to B28
$wsFunciones2 = $wb.Worksheets.Item(3)
$wsFunciones2.Range("B28").Select()

# Sheet 4 "Calificaciones ponderada": selection moves to E31 (no longer the active tab)
$wsCalificaciones = $wb.Worksheets.Item(4)
$wsCalificaciones.Range("E31").Select()

# Restore the active sheet/selection to sheet 1, since that is the new active tab
$wsInversion.Activate()
$wsInversion.Range("K12").Select()
